$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L, mirroring column K's 2020 data (same value + formatting)
# Row 3 (year header): 2020
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L3").Value = 2020

# Row 4 (data value): 6.18
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L4").Value = 6.18

# Update the saved selection/active cell to M12, as recorded in the view state
$ws.Range("M12").Select()
